$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete data rows 2-7, keeping only the header row
$ws.Range("A2:D7").EntireRow.Delete()

# Update header text for column B
$ws.Range("B1").Value = "Valor aproximado"

# Add new header column E
$ws.Range("E1").Value = "Tipo bolo"
